$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1 in the package, first tab)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5153
$ws1.Range("F7").Value = 782

# Update "全部类型" sheet (sheet4 in the package, fourth tab)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5153
$ws4.Range("F7").Value = 782
